# Apply updated TPM-derived statistics to the NATMI ligand-receptor pair sheet.
# (Re-running the pipeline with new TPM values changed only the numeric result
#  columns; the Sending/Target cluster and Ligand/Receptor labels are unchanged.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.9352601111131627
$ws.Range("J2").Value = 0.9352601111131628
$ws.Range("M2").Value = 34.08558233333333
$ws.Range("N2").Value = 102.256747
$ws.Range("O2").Value = 0.3376420874206352
$ws.Range("P2").Value = 0.3376420874206352
$ws.Range("Q2").Value = 10.93489341160966
$ws.Range("R2").Value = 98.41404070448699
$ws.Range("S2").Value = 0.3157831761975035
$ws.Range("T2").Value = 0.3157831761975035

# Row 3
$ws.Range("I3").Value = 0.9352601111131627
$ws.Range("J3").Value = 0.9352601111131628
$ws.Range("N3").Value = 90.767769
$ws.Range("O3").Value = 0.2997065709089496
$ws.Range("P3").Value = 0.2997065709089496
$ws.Range("Q3").Value = 9.706311889861
$ws.Range("R3").Value = 87.356807008749
$ws.Range("S3").Value = 0.2803036008096492
$ws.Range("T3").Value = 0.2803036008096492

# Row 4
$ws.Range("I4").Value = 0.9352601111131627
$ws.Range("J4").Value = 0.9352601111131628
$ws.Range("M4").Value = 8.754337666666666
$ws.Range("N4").Value = 26.263013
$ws.Range("O4").Value = 0.08671797990283496
$ws.Range("P4").Value = 0.08671797990283496
$ws.Range("Q4").Value = 2.808452803830333
$ws.Range("R4").Value = 25.276075234473
$ws.Range("S4").Value = 0.08110386751943444
$ws.Range("T4").Value = 0.08110386751943444

# Row 5
$ws.Range("I5").Value = 0.9352601111131627
$ws.Range("J5").Value = 0.9352601111131628
$ws.Range("M5").Value = 27.85597433333334
$ws.Range("N5").Value = 83.56792300000001
$ws.Range("O5").Value = 0.2759333617675802
$ws.Range("P5").Value = 0.2759333617675802
$ws.Range("Q5").Value = 8.936391557953668
$ws.Range("R5").Value = 80.42752402158301
$ws.Range("S5").Value = 0.2580694665865755
$ws.Range("T5").Value = 0.2580694665865755

# Row 6
$ws.Range("G6").Value = 0.02220666666666667
$ws.Range("H6").Value = 0.06662
$ws.Range("I6").Value = 0.06473988888683736
$ws.Range("J6").Value = 0.06473988888683736
$ws.Range("M6").Value = 34.08558233333333
$ws.Range("N6").Value = 102.256747
$ws.Range("O6").Value = 0.3376420874206352
$ws.Range("P6").Value = 0.3376420874206352
$ws.Range("Q6").Value = 0.7569271650155555
$ws.Range("R6").Value = 6.812344485139999
$ws.Range("S6").Value = 0.02185891122313175
$ws.Range("T6").Value = 0.02185891122313175

# Row 7
$ws.Range("G7").Value = 0.02220666666666667
$ws.Range("H7").Value = 0.06662
$ws.Range("I7").Value = 0.06473988888683736
$ws.Range("J7").Value = 0.06473988888683736
$ws.Range("N7").Value = 90.767769
$ws.Range("O7").Value = 0.2997065709089496
$ws.Range("P7").Value = 0.2997065709089496
$ws.Range("Q7").Value = 0.6718831967533333
$ws.Range("R7").Value = 6.04694877078
$ws.Range("S7").Value = 0.01940297009930044
$ws.Range("T7").Value = 0.01940297009930044

# Row 8
$ws.Range("G8").Value = 0.02220666666666667
$ws.Range("H8").Value = 0.06662
$ws.Range("I8").Value = 0.06473988888683736
$ws.Range("J8").Value = 0.06473988888683736
$ws.Range("M8").Value = 8.754337666666666
$ws.Range("N8").Value = 26.263013
$ws.Range("O8").Value = 0.08671797990283496
$ws.Range("P8").Value = 0.08671797990283496
$ws.Range("Q8").Value = 0.1944046584511111
$ws.Range("R8").Value = 1.74964192606
$ws.Range("S8").Value = 0.005614112383400531
$ws.Range("T8").Value = 0.005614112383400531

# Row 9
$ws.Range("G9").Value = 0.02220666666666667
$ws.Range("H9").Value = 0.06662
$ws.Range("I9").Value = 0.06473988888683736
$ws.Range("J9").Value = 0.06473988888683736
$ws.Range("M9").Value = 27.85597433333334
$ws.Range("N9").Value = 83.56792300000001
$ws.Range("O9").Value = 0.2759333617675802
$ws.Range("P9").Value = 0.2759333617675802
$ws.Range("Q9").Value = 0.6185883366955556
$ws.Range("R9").Value = 5.56729503026
$ws.Range("S9").Value = 0.01786389518100464
$ws.Range("T9").Value = 0.01786389518100464
